$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Personalnummer value from M100002 to M100001
$ws.Range("B2").Value = "M100001"

# Move the active selection from B4 to B3
$ws.Range("B3").Select()
